$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the date style for the new row 5 birthday-date cell before putting a value in,
# mirroring the other date cells (C2:C4). Use copy/paste-special so the existing
# numFmtId=14 style (s=1) is reused instead of a new custom number format being created.
$ws.Cells.Item(4, 3).Copy() | Out-Null
$ws.Cells.Item(5, 3).PasteSpecial(-4122) | Out-Null
# Email column style (hyperlink look) for the new/filled rows, same as D2/D3.
$ws.Cells.Item(4, 4).Style = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(5, 4).Style = $ws.Cells.Item(3, 4).Style

# Order matters for shared-string table layout, so we follow the same
# sequence the original author used when typing values in:
# D5 email text, B4 name, D4 email text, B5 name.
$ws.Cells.Item(5, 4).Value = "deendayal555kumawat@gmail.com"
$ws.Cells.Item(4, 2).Value = "Iron Man"
$ws.Cells.Item(4, 4).Value = "user221user@gmail.com"
$ws.Cells.Item(5, 2).Value = "Captain America"

# Hyperlinks are added in D4, D5 order so relationship ids come out as rId3, rId4.
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:user221user@gmail.com", [Type]::Missing, [Type]::Missing, "user221user@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:deendayal555kumawat@gmail.com", [Type]::Missing, [Type]::Missing, "deendayal555kumawat@gmail.com")

# Re-apply the plain (non-hyperlink-styled) look to match the rest of the sheet.
$ws.Cells.Item(4, 4).Style = $ws.Cells.Item(3, 4).Style
$ws.Cells.Item(5, 4).Style = $ws.Cells.Item(3, 4).Style

# Remaining values for row 4
$ws.Cells.Item(4, 5).Value = 7648648887
$ws.Cells.Item(4, 6).Value = "fknvbckjbnkcjFHFHGFkhkjhkjhkj76576jdfkjdh"

# Remaining values for row 5
$ws.Cells.Item(5, 1).Value = 1906
$ws.Cells.Item(5, 3).Value = 36624
$ws.Cells.Item(5, 5).Value = 87674876845
$ws.Cells.Item(5, 6).Value = "dfjdjJKGJJBKJKjlkdhfjkgdhkj983457843dbffdj"

$ws.Range("F5").Select() | Out-Null
